$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)      # "总计"
$q2 = $wb.Worksheets.Item(2)           # currently "2022-Q2" - will become "2022-Q4"

# --- Duplicate the "2022-Q2" sheet so the original fund-holdings data is
# preserved on its own tab, positioned right after the renamed sheet. ---
$q2.Copy($null, $q2)
$q2copy = $wb.Worksheets.Item(3)

# Rename sheets: the original tab becomes "2022-Q4" (new data), the
# duplicate keeps the old "2022-Q2" data under its original name.
$q2.Name = "2022-Q4"
$q2copy.Name = "2022-Q2"
$q4 = $q2

# --- Update the "总计" (summary) sheet ---
# Row 2 now reports the new quarter's totals.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("D2").Value = 0.62

# Row 3 is a new row holding what used to be the latest entry.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0

# --- Update the "2022-Q4" sheet with the new fund table ---
# Bring over the bold/bordered header style used on the summary sheet.
$summary.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)

# Write the new fund row. Text-like numeric values (fund code, sizes,
# percentages) must stay text, so build them via a text formula in a
# scratch cell, then paste-values into place (avoids leading-zero loss
# and avoids leaving a stray quote-prefixed style behind).
function Set-TextValue($range, [string]$text) {
    $scratch = $q4.Range("Z100")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue $q4.Range("B2") "010695"
Set-TextValue $q4.Range("C2") "华夏磐益一年定期开放混合"
Set-TextValue $q4.Range("D2") "16.03"
Set-TextValue $q4.Range("E2") "98.69"
Set-TextValue $q4.Range("F2") "3.88"
Set-TextValue $q4.Range("G2") "0.6220"
$q4.Range("H2").Value = 9

# Match the page margins used on the summary sheet (values are in points).
$q4.PageSetup.LeftMargin = 0.75 * 72
$q4.PageSetup.RightMargin = 0.75 * 72
$q4.PageSetup.TopMargin = 1 * 72
$q4.PageSetup.BottomMargin = 1 * 72
$q4.PageSetup.HeaderMargin = 0.5 * 72
$q4.PageSetup.FooterMargin = 0.5 * 72
